# Add a new (3rd) slide using the Blank layout and place a transparent
# (no fill / no line) rectangle on it — used as a positioning/background
# placeholder for an icon image.

$p = $ppt.ActivePresentation

# ppLayoutBlank = 12 -> the slide layout named "Blank" (slideLayout7.xml).
$s = $p.Slides.Add(3, 12)

# The target shape ends up with id=4 / name "Rectangle 3": shape ids/names
# are assigned sequentially per-slide (id 1 is the root group), so add and
# discard two throw-away rectangles first to advance the counters the same
# way the authoring session that produced the target deck did.
$tmp1 = $s.Shapes.AddShape(1, 0, 0, 10, 10)
$tmp1.Delete()
$tmp2 = $s.Shapes.AddShape(1, 0, 0, 10, 10)
$tmp2.Delete()

# EMU -> point conversion (1 pt = 12700 EMU) since Shapes.AddShape takes
# points for its position/size arguments.
$shp = $s.Shapes.AddShape(1, (3580760 / 12700.0), (1759644 / 12700.0), (1129553 / 12700.0), (1091132 / 12700.0))

# Transparent background: no fill, no outline.
$shp.Fill.Visible = 0
$shp.Line.Visible = 0

# Center the (empty) text, matching the authored shape's paragraph/body
# properties.
$shp.TextFrame.VerticalAnchor = 3
$shp.TextFrame.TextRange.ParagraphFormat.Alignment = 2
$shp.TextFrame.TextRange.LanguageID = "en-AU"

Write-Host "Added slide $($s.SlideIndex) with shape '$($shp.Name)' (id $($shp.Id))"
